$d = $word.ActiveDocument

# Locate the target paragraph (the one starting with "Модель линейной регресси...")
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Модель линейной*") {
        $target = $p
    }
}
if ($target -eq $null) {
    throw "Target paragraph not found"
}

$full = $target.Range
# Range covering the paragraph's content but excluding the trailing paragraph mark
# (the empty _GoBack bookmark normally sits right before the mark).
$content = $d.Range($full.Start, $full.End - 1)

# Phase 1: rewrite everything up to (but not including) the final run/sentence.
# This intentionally omits the trailing bookmark - Word re-synthesizes a fresh
# _GoBack around whatever range we touch, and phase 2 below re-targets it
# precisely to a zero-width location at the paragraph's end.
$xmlFrag1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>М</w:t></w:r><w:r><w:t xml:space="preserve">одель линейной </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>регресси</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> с L1-регуляризацией</w:t></w:r><w:r><w:t xml:space="preserve"> показала себя лучше всего. </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>L</w:t></w:r><w:r w:rsidRPr="00EA02C6"><w:t xml:space="preserve">1 </w:t></w:r><w:r><w:t xml:space="preserve">регуляризация </w:t></w:r><w:r><w:t>нужна для</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>очистки данных</w:t></w:r><w:r><w:t xml:space="preserve"> от</w:t></w:r><w:r><w:t xml:space="preserve"> шума</w:t></w:r><w:r><w:t>Z</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$content.InsertXML($xmlFrag1)

# Remove the bookmark Word auto-created around the whole freshly inserted range;
# it will be re-created precisely (zero-width) in phase 2.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Phase 2: locate the sentinel "Z" we appended above, which now sits at the very
# end of the paragraph's content. Replacing just that one character lets the new
# content (the real last run's text) land exactly at the paragraph's end, with a
# zero-width _GoBack bookmark synthesized right after it.
$target2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*шумаZ*") {
        $target2 = $p
    }
}
if ($target2 -eq $null) {
    throw "Sentinel paragraph not found"
}
$p2End = $target2.Range.End - 1
$sentinel = $d.Range($p2End - 1, $p2End)

$xmlFrag2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>. Из этого можно сделать вывод что в наших данных было много шума.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$sentinel.InsertXML($xmlFrag2)

Write-Output "Done"
